$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("B2").Style

$cell = $ws.Range("D2")
$cell.Value = "'29.049.53"
$cell.Style = $plainStyle
$cell = $ws.Range("E2")
$cell.Value = "'  -0.63%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D3")
$cell.Value = "'1.831.63"
$cell.Style = $plainStyle
$cell = $ws.Range("E3")
$cell.Value = "'  -0.70%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D4")
$cell.Value = "'0.9993"
$cell.Style = $plainStyle
$cell = $ws.Range("E4")
$cell.Value = "'  +0.04%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D5")
$cell.Value = "'241.72"
$cell.Style = $plainStyle
$cell = $ws.Range("E5")
$cell.Value = "'  +0.45%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D6")
$cell.Value = "'0.6550"
$cell.Style = $plainStyle
$cell = $ws.Range("E6")
$cell.Value = "'  -2.67%  "
$cell.Style = $plainStyle
$cell = $ws.Range("E7")
$cell.Value = "'  +0.04%  "
$cell.Style = $plainStyle
$cell = $ws.Range("E8")
$cell.Value = "'  +5.82%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D9")
$cell.Value = "'0.2936"
$cell.Style = $plainStyle
$cell = $ws.Range("E9")
$cell.Value = "'  -0.58%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D10")
$cell.Value = "'0.07337"
$cell.Style = $plainStyle
$cell = $ws.Range("E10")
$cell.Value = "'  -1.23%  "
$cell.Style = $plainStyle
$cell = $ws.Range("E11")
$cell.Value = "'  +0.47%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D12")
$cell.Value = "'0.07672"
$cell.Style = $plainStyle
$cell = $ws.Range("E12")
$cell.Value = "'  -0.56%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D13")
$cell.Value = "'1.841.26"
$cell.Style = $plainStyle
$cell = $ws.Range("E13")
$cell.Value = "'  +0.27%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D14")
$cell.Value = "'4.981"
$cell.Style = $plainStyle
$cell = $ws.Range("E14")
$cell.Value = "'  -0.54%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D15")
$cell.Value = "'0.6675"
$cell.Style = $plainStyle
$cell = $ws.Range("E15")
$cell.Value = "'  -0.70%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D16")
$cell.Value = "'81.70"
$cell.Style = $plainStyle
$cell = $ws.Range("E16")
$cell.Value = "'  -5.20%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D17")
$cell.Value = "'6.107"
$cell.Style = $plainStyle
$cell = $ws.Range("E17")
$cell.Value = "'  -0.23%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D18")
$cell.Value = "'0.000008669"
$cell.Style = $plainStyle
$cell = $ws.Range("E18")
$cell.Value = "'  +4.26%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D19")
$cell.Value = "'29.042.69"
$cell.Style = $plainStyle
$cell = $ws.Range("E19")
$cell.Value = "'  -0.39%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D20")
$cell.Value = "'2.089.29"
$cell.Style = $plainStyle
$cell = $ws.Range("E20")
$cell.Value = "'  +1.07%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D21")
$cell.Value = "'12.44"
$cell.Style = $plainStyle
$cell = $ws.Range("E21")
$cell.Value = "'  -0.61%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D22")
$cell.Value = "'223.95"
$cell.Style = $plainStyle
$cell = $ws.Range("E22")
$cell.Value = "'  -2.07%  "
$cell.Style = $plainStyle
$cell = $ws.Range("E23")
$cell.Value = "'  -0.04%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D24")
$cell.Value = "'7.135"
$cell.Style = $plainStyle
$cell = $ws.Range("E24")
$cell.Value = "'  -0.72%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D26")
$cell.Value = "'157.56"
$cell.Style = $plainStyle
$cell = $ws.Range("E26")
$cell.Value = "'  -1.89%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D27")
$cell.Value = "'8.495"
$cell.Style = $plainStyle
$cell = $ws.Range("E27")
$cell.Value = "'  -2.25%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D28")
$cell.Value = "'0.1380"
$cell.Style = $plainStyle
$cell = $ws.Range("E28")
$cell.Value = "'  -1.73%  "
$cell.Style = $plainStyle
$cell = $ws.Range("E29")
$cell.Value = "'  -0.47%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D30")
$cell.Value = "'1.508"
$cell.Style = $plainStyle
$cell = $ws.Range("E30")
$cell.Value = "'  +0.04%  "
$cell.Style = $plainStyle
$cell = $ws.Range("E31")
$cell.Value = "'  -1.75%  "
$cell.Style = $plainStyle
$cell = $ws.Range("B32")
$cell.Value = "'InternetComputer(DFINITY)"
$cell.Style = $plainStyle
$cell = $ws.Range("C32")
$cell.Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell.Style = $plainStyle
$cell = $ws.Range("D32")
$cell.Value = "'4.012"
$cell.Style = $plainStyle
$cell = $ws.Range("E32")
$cell.Value = "'  -1.41%  "
$cell.Style = $plainStyle
$cell = $ws.Range("B33")
$cell.Value = "'Toncoin"
$cell.Style = $plainStyle
$cell = $ws.Range("C33")
$cell.Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell.Style = $plainStyle
$cell = $ws.Range("D33")
$cell.Value = "'1.200"
$cell.Style = $plainStyle
$cell = $ws.Range("E33")
$cell.Value = "'  +0.80%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D34")
$cell.Value = "'0.05350"
$cell.Style = $plainStyle
$cell = $ws.Range("E34")
$cell.Value = "'  +0.84%  "
$cell.Style = $plainStyle
$cell = $ws.Range("E35")
$cell.Value = "'  -1.96%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D36")
$cell.Value = "'0.7425"
$cell.Style = $plainStyle
$cell = $ws.Range("E36")
$cell.Value = "'  -2.03%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D37")
$cell.Value = "'1.160"
$cell.Style = $plainStyle
$cell = $ws.Range("E37")
$cell.Value = "'  +2.08%  "
$cell.Style = $plainStyle
$cell = $ws.Range("E38")
$cell.Value = "'  -0.91%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D39")
$cell.Value = "'1.297.84"
$cell.Style = $plainStyle
$cell = $ws.Range("E39")
$cell.Value = "'  -2.20%  "
$cell.Style = $plainStyle
$cell = $ws.Range("B40")
$cell.Value = "'VeChain"
$cell.Style = $plainStyle
$cell = $ws.Range("C40")
$cell.Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell.Style = $plainStyle
$cell = $ws.Range("D40")
$cell.Value = "'0.01786"
$cell.Style = $plainStyle
$cell = $ws.Range("E40")
$cell.Value = "'  -1.00%  "
$cell.Style = $plainStyle
$cell = $ws.Range("B41")
$cell.Value = "'MXToken"
$cell.Style = $plainStyle
$cell = $ws.Range("C41")
$cell.Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell.Style = $plainStyle
$cell = $ws.Range("D41")
$cell.Value = "'2.750"
$cell.Style = $plainStyle
$cell = $ws.Range("E41")
$cell.Value = "'  +0.75%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D42")
$cell.Value = "'6.328"
$cell.Style = $plainStyle
$cell = $ws.Range("E42")
$cell.Value = "'  +6.28%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D43")
$cell.Value = "'0.9025"
$cell.Style = $plainStyle
$cell = $ws.Range("E43")
$cell.Value = "'  -1.80%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D44")
$cell.Value = "'0.9997"
$cell.Style = $plainStyle
$cell = $ws.Range("E44")
$cell.Value = "'  -0.17%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D45")
$cell.Value = "'103.24"
$cell.Style = $plainStyle
$cell = $ws.Range("E45")
$cell.Value = "'  -0.17%  "
$cell.Style = $plainStyle
$cell = $ws.Range("B46")
$cell.Value = "'XinFinNetwork"
$cell.Style = $plainStyle
$cell = $ws.Range("C46")
$cell.Value = "'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$cell.Style = $plainStyle
$cell = $ws.Range("D46")
$cell.Value = "'0.07975"
$cell.Style = $plainStyle
$cell = $ws.Range("E46")
$cell.Value = "'  -0.86%  "
$cell.Style = $plainStyle
$cell = $ws.Range("B47")
$cell.Value = "'RocketPoolETH"
$cell.Style = $plainStyle
$cell = $ws.Range("C47")
$cell.Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$cell.Style = $plainStyle
$cell = $ws.Range("D47")
$cell.Value = "'1.987.92"
$cell.Style = $plainStyle
$cell = $ws.Range("E47")
$cell.Value = "'  +0.88%  "
$cell.Style = $plainStyle
$cell = $ws.Range("B48")
$cell.Value = "'Aave"
$cell.Style = $plainStyle
$cell = $ws.Range("C48")
$cell.Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell.Style = $plainStyle
$cell = $ws.Range("D48")
$cell.Value = "'64.52"
$cell.Style = $plainStyle
$cell = $ws.Range("E48")
$cell.Value = "'  +0.99%  "
$cell.Style = $plainStyle
$cell = $ws.Range("B49")
$cell.Value = "'Mantle"
$cell.Style = $plainStyle
$cell = $ws.Range("C49")
$cell.Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell.Style = $plainStyle
$cell = $ws.Range("D49")
$cell.Value = "'0.5137"
$cell.Style = $plainStyle
$cell = $ws.Range("E49")
$cell.Value = "'  -0.47%  "
$cell.Style = $plainStyle
$cell = $ws.Range("B50")
$cell.Value = "'BabyDogeCoin"
$cell.Style = $plainStyle
$cell = $ws.Range("C50")
$cell.Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$cell.Style = $plainStyle
$cell = $ws.Range("D50")
$cell.Value = "'0.00000000121"
$cell.Style = $plainStyle
$cell = $ws.Range("E50")
$cell.Value = "'  -0.15%  "
$cell.Style = $plainStyle
$cell = $ws.Range("D51")
$cell.Value = "'1.739"
$cell.Style = $plainStyle
$cell = $ws.Range("E51")
$cell.Value = "'  -1.96%  "
$cell.Style = $plainStyle
